# Apply the "I0"/"IF" columns (I and J) edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 = "I0" and J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting (bold font + border + centered alignment)
# already used by the other header cells (copy format from H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-59: add numeric values for I and J ---
$rowData = @(
    @{Row=2; I=7; J=7},
    @{Row=3; I=7; J=7},
    @{Row=4; I=5; J=5},
    @{Row=5; I=8; J=8},
    @{Row=6; I=5; J=5},
    @{Row=7; I=7; J=7},
    @{Row=8; I=2; J=2},
    @{Row=9; I=7; J=7},
    @{Row=10; I=8; J=8},
    @{Row=11; I=5; J=5},
    @{Row=12; I=7; J=8},
    @{Row=13; I=8; J=8},
    @{Row=14; I=6; J=6},
    @{Row=15; I=8; J=8},
    @{Row=16; I=8; J=8},
    @{Row=17; I=7; J=7},
    @{Row=18; I=9; J=9},
    @{Row=19; I=5; J=5},
    @{Row=20; I=5; J=5},
    @{Row=21; I=5; J=6},
    @{Row=22; I=4; J=4},
    @{Row=23; I=1; J=1},
    @{Row=24; I=6; J=6},
    @{Row=25; I=1; J=2},
    @{Row=26; I=6; J=6},
    @{Row=27; I=8; J=8},
    @{Row=28; I=1; J=2},
    @{Row=29; I=6; J=6},
    @{Row=30; I=1; J=2},
    @{Row=31; I=8; J=8},
    @{Row=32; I=1; J=2},
    @{Row=33; I=1; J=2},
    @{Row=34; I=2; J=3},
    @{Row=35; I=6; J=7},
    @{Row=36; I=7; J=8},
    @{Row=37; I=6; J=6},
    @{Row=38; I=9; J=9},
    @{Row=39; I=8; J=8},
    @{Row=40; I=7; J=7},
    @{Row=41; I=8; J=8},
    @{Row=42; I=4; J=4},
    @{Row=43; I=7; J=8},
    @{Row=44; I=7; J=8},
    @{Row=45; I=9; J=9},
    @{Row=46; I=6; J=7},
    @{Row=47; I=6; J=6},
    @{Row=48; I=8; J=8},
    @{Row=49; I=5; J=5},
    @{Row=50; I=7; J=8},
    @{Row=51; I=6; J=7},
    @{Row=52; I=6; J=6},
    @{Row=53; I=9; J=9},
    @{Row=54; I=10; J=10},
    @{Row=55; I=7; J=7},
    @{Row=56; I=7; J=8},
    @{Row=57; I=7; J=9},
    @{Row=58; I=7; J=8},
    @{Row=59; I=3; J=3}
)

foreach ($item in $rowData) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
